# "excel sums to 1440 on every column"
# Correct a handful of per-country time-use figures (column B, "Total") so
# each row's minutes sum to 1440 (24h). A few of the corrections are left
# as formulas mirroring the authored edit (offset arithmetic from the old
# value), the rest are plain corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (Japan): 325.7113724006864 -> 301.711372400686 (via formula, -24)
$ws.Range("B15").Formula = "=325.711372400686-24"

# Row 15, column L (Japan, "Sports"): 5 -> 0
$ws.Range("L15").Value = 0

# Row 20 (Mexico): 302.33356400000002 -> 304.33356400000002
$ws.Range("B20").Value = 304.33356400000002

# Row 22 (New Zealand): 240.99999999999997 -> 239 (via formula, 241-2)
$ws.Range("B22").Formula = "=241-2"

# Row 28 (Sweden): 262.18713699086243 -> 268.18713699086197
$ws.Range("B28").Value = 268.18713699086197

# Row 29 (Turkey): 217.02784000000003 -> 217.02784
$ws.Range("B29").Value = 217.02784

# Row 30 (UK): 235.49323385367586 -> 234.493233853676
$ws.Range("B30").Value = 234.493233853676

# Row 34 (South Africa): 188.80248 -> 201.80248 (via formula, +13)
$ws.Range("B34").Formula = "=188.80248+13"

# Update the saved selection / active cell shown when the workbook is reopened
$ws.Range("T17").Select() | Out-Null
